$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2185.5908
$ws.Range("I19").Value = 2142.75
$ws.Range("J19").Value = 2299.8333
$ws.Range("K19").Value = 2142.75
$ws.Range("L19").Value = 2299.8333
$ws.Range("M19").Value = -1967.75
$ws.Range("N19").Value = -2649.8333
$ws.Range("H55").Value = 416.6
$ws.Range("I55").Value = 361
$ws.Range("K55").Value = 361
$ws.Range("M55").Value = -147
$ws.Range("H74").Value = 4460.3335
$ws.Range("I74").Value = 4268
$ws.Range("J74").Value = 5999
$ws.Range("K74").Value = 4268
$ws.Range("L74").Value = 5999
$ws.Range("M74").Value = -3332
$ws.Range("N74").Value = -7871
$ws.Range("H77").Value = 4460.3335
$ws.Range("I77").Value = 4268
$ws.Range("J77").Value = 5999
$ws.Range("K77").Value = 21340
$ws.Range("L77").Value = 29995
$ws.Range("M77").Value = -16660
$ws.Range("N77").Value = -39355
$ws.Range("H112").Value = 3812.6667
$ws.Range("J112").Value = 3812.6667
$ws.Range("L112").Value = 11438.0001
$ws.Range("N112").Value = -13654.0001
$ws.Range("H116").Value = 24642426
$ws.Range("I116").Value = 33337856
$ws.Range("J116").Value = 5374.5
$ws.Range("K116").Value = 33337856
$ws.Range("L116").Value = 5374.5
$ws.Range("M116").Value = -33334414
$ws.Range("N116").Value = -12258.5
$ws.Range("H133").Value = 59250
$ws.Range("J133").Value = 59250
$ws.Range("L133").Value = 59250
$ws.Range("N133").Value = -69370
$ws.Range("H138").Value = 5538.74
$ws.Range("J138").Value = 7143.2456
$ws.Range("L138").Value = 21429.7368
$ws.Range("N138").Value = -31709.7368

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 865.1923
$ws.Range("I2").Value = 878.04346
$ws.Range("K2").Value = 878.04346
$ws.Range("M2").Value = -765.04346
$ws.Range("H45").Value = 2310.1177
$ws.Range("I45").Value = 1670.1818
$ws.Range("K45").Value = 1670.1818
$ws.Range("M45").Value = -1293.1818
$ws.Range("H74").Value = 1256.7142
$ws.Range("I74").Value = 1221.7778
$ws.Range("K74").Value = 1221.7778
$ws.Range("M74").Value = -347.7778000000001
$ws.Range("H77").Value = 1256.7142
$ws.Range("I77").Value = 1221.7778
$ws.Range("K77").Value = 6108.889
$ws.Range("M77").Value = -1740.889
$ws.Range("H109").Value = 57259
$ws.Range("J109").Value = 57259
$ws.Range("L109").Value = 57259
$ws.Range("N109").Value = -60033
$ws.Range("H116").Value = 865.1923
$ws.Range("I116").Value = 878.04346
$ws.Range("K116").Value = 878.04346
$ws.Range("M116").Value = 1415.95654
$ws.Range("H122").Value = 3443.7778
$ws.Range("I122").Value = 3443.7778
$ws.Range("K122").Value = 10331.3334
$ws.Range("M122").Value = -7881.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 865.1923
$ws.Range("I3").Value = 878.04346
$ws.Range("K3").Value = 878.04346
$ws.Range("M3").Value = -764.04346

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 45296.5
$ws.Range("J52").Value = 46969.75
$ws.Range("L52").Value = 46969.75
$ws.Range("N52").Value = -47557.75
$ws.Range("H107").Value = 604.6667
$ws.Range("J107").Value = 761.2727
$ws.Range("L107").Value = 761.2727
$ws.Range("N107").Value = -4601.2727
$ws.Range("H134").Value = 3016.5557
$ws.Range("I134").Value = 2003.6875
$ws.Range("K134").Value = 6011.0625
$ws.Range("M134").Value = -3476.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 5069.4287
$ws.Range("I69").Value = 4395
$ws.Range("J69").Value = 5444.1113
$ws.Range("K69").Value = 13185
$ws.Range("L69").Value = 16332.3339
$ws.Range("M69").Value = -12374
$ws.Range("N69").Value = -17954.3339
$ws.Range("H72").Value = 5069.4287
$ws.Range("I72").Value = 4395
$ws.Range("J72").Value = 5444.1113
$ws.Range("K72").Value = 39555
$ws.Range("L72").Value = 48997.00169999999
$ws.Range("M72").Value = -35499
$ws.Range("N72").Value = -57109.00169999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2352.6155
$ws.Range("I97").Value = 1413
$ws.Range("J97").Value = 3158
$ws.Range("K97").Value = 1413
$ws.Range("L97").Value = 3158
$ws.Range("M97").Value = -917
$ws.Range("N97").Value = -4150
$ws.Range("H114").Value = 91233
$ws.Range("J114").Value = 91233
$ws.Range("L114").Value = 91233
$ws.Range("N114").Value = -99911
$ws.Range("H132").Value = 119523.47
$ws.Range("I132").Value = 155190.39
$ws.Range("J132").Value = 3606
$ws.Range("K132").Value = 465571.17
$ws.Range("L132").Value = 10818
$ws.Range("M132").Value = -463041.17
$ws.Range("N132").Value = -15878
$ws.Range("H136").Value = 32122.615
$ws.Range("J136").Value = 32122.615
$ws.Range("L136").Value = 96367.845
$ws.Range("N136").Value = -101467.845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5916.364
$ws.Range("I61").Value = 6276.222
$ws.Range("K61").Value = 6276.222
$ws.Range("M61").Value = -6074.222
$ws.Range("H113").Value = 5916.364
$ws.Range("I113").Value = 6276.222
$ws.Range("K113").Value = 6276.222
$ws.Range("M113").Value = -4106.222
$ws.Range("H122").Value = 6607.326
$ws.Range("I122").Value = 4006.2068
$ws.Range("K122").Value = 12018.6204
$ws.Range("M122").Value = -9568.6204
$ws.Range("H132").Value = 4060.8987
$ws.Range("I132").Value = 3297.5
$ws.Range("J132").Value = 5020.6
$ws.Range("K132").Value = 9892.5
$ws.Range("L132").Value = 15061.8
$ws.Range("M132").Value = -7362.5
$ws.Range("N132").Value = -20121.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8497.786
$ws.Range("I81").Value = 6066.5
$ws.Range("J81").Value = 10321.25
$ws.Range("K81").Value = 12133
$ws.Range("L81").Value = 20642.5
$ws.Range("M81").Value = -11072
$ws.Range("N81").Value = -22764.5
$ws.Range("H84").Value = 8497.786
$ws.Range("I84").Value = 6066.5
$ws.Range("J84").Value = 10321.25
$ws.Range("K84").Value = 60665
$ws.Range("L84").Value = 103212.5
$ws.Range("M84").Value = -55361
$ws.Range("N84").Value = -113820.5
$ws.Range("H96").Value = 1391.5
$ws.Range("J96").Value = 1410
$ws.Range("L96").Value = 1410
$ws.Range("N96").Value = -4156
$ws.Range("H113").Value = 361.9091
$ws.Range("I113").Value = 391.875
$ws.Range("J113").Value = 282
$ws.Range("K113").Value = 1175.625
$ws.Range("L113").Value = 846
$ws.Range("M113").Value = 994.375
$ws.Range("N113").Value = -5186
$ws.Range("H122").Value = 3157.6843
$ws.Range("I122").Value = 2099.7273
$ws.Range("J122").Value = 4612.375
$ws.Range("K122").Value = 6299.1819
$ws.Range("L122").Value = 13837.125
$ws.Range("M122").Value = -3849.1819
$ws.Range("N122").Value = -18737.125
$ws.Range("H136").Value = 3109.2
$ws.Range("I136").Value = 2237.0938
$ws.Range("K136").Value = 6711.2814
$ws.Range("M136").Value = -4161.2814
